$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text values: Excel's literal-input type inference leaves these as
# text already (they aren't valid numbers/dates), so a normal .Value
# assignment is enough and keeps the default style (no NumberFormat churn).
function Set-PlainText($addr, $value) {
    $ws.Range($addr).Value = $value
}

# Numeric-looking / date-looking strings: a direct .Value assignment would
# get auto-coerced by Excel into a real number or date serial (and a date
# also drags in a new date-formatted style). Routing the literal through a
# text formula and then collapsing it back to a value via Paste Special
# (values only) keeps the original text as a genuine shared string without
# adding any NumberFormat/style to the cell.
function Set-TextLiteral($addr, $value) {
    $r = $ws.Range($addr)
    $escaped = $value.Replace('"', '""')
    $r.Formula = '="' + $escaped + '"'
    $r.Copy()
    $r.PasteSpecial(-4163)
}

Set-PlainText   "C1"  "FREECODE"
Set-TextLiteral "D1"  "10-06-2020"
Set-PlainText   "E1"  "FA"
Set-PlainText   "F1"  "https://www.youtube.com/watch?v=9sWEecNUW-o"
Set-PlainText   "G1"  "Code your own YouTube app: YouTube API + HTML + CSS + JavaScript (full tutorial)"
Set-PlainText   "H1"  "1-TestFile"
Set-PlainText   "I1"  "PT1H7M35S"
Set-PlainText   "J1"  "freeCodeCamp.org"
Set-TextLiteral "Q1"  "2060000"
Set-PlainText   "AD1" "2014-12-16T21:18:48Z"
Set-TextLiteral "AE1" "1114"
Set-TextLiteral "AH1" "74760"
Set-TextLiteral "AI1" "1297"
Set-TextLiteral "AJ1" "46"
Set-TextLiteral "AX1" "96"
